$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.920.31"
$ws.Range("E2").Value = "  +4.32%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.079.13"
$ws.Range("E3").Value = "  +2.83%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.82"
$ws.Range("E5").Value = "  +2.92%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.27"
$ws.Range("E6").Value = "  +2.32%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.069.37"
$ws.Range("E8").Value = "  +2.93%  "

$ws.Range("E9").Value = "  +1.09%  "

$ws.Range("E10").Value = "  +5.05%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.78"
$ws.Range("E11").Value = "  +11.61%  "

$ws.Range("E12").Value = "  +2.18%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000240"
$ws.Range("E13").Value = "  +4.16%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.33"
$ws.Range("E14").Value = "  +4.42%  "

$ws.Range("E15").Value = "  +0.13%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.588.07"
$ws.Range("E16").Value = "  +3.15%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.25"
$ws.Range("E17").Value = "  -0.72%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.077.38"
$ws.Range("E18").Value = "  +3.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "61.808.04"
$ws.Range("E19").Value = "  +4.45%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "447.55"
$ws.Range("E20").Value = "  +4.14%  "

$ws.Range("E21").Value = "  +1.97%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.731"
$ws.Range("E22").Value = "  +1.70%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.45"
$ws.Range("E23").Value = "  +4.41%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.82"
$ws.Range("E24").Value = "  +3.11%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.83"
$ws.Range("E25").Value = "  +0.93%  "

$ws.Range("E26").Value = "  +0.31%  "

$ws.Range("E27").Value = "  +4.35%  "

$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.28"
$ws.Range("E28").Value = "  +7.07%  "

$ws.Range("B29").Value = "FirstDigitalUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.67"
$ws.Range("E30").Value = "  +4.77%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.78"
$ws.Range("E31").Value = "  +11.02%  "

$ws.Range("E32").Value = "  +12.29%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.82"
$ws.Range("E33").Value = "  +4.06%  "

$ws.Range("E34").Value = "  +4.97%  "

$ws.Range("E35").Value = "  +2.73%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.05"
$ws.Range("E36").Value = "  +2.60%  "

$ws.Range("E37").Value = "  +4.81%  "

$ws.Range("E38").Value = "  +1.89%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.00"
$ws.Range("E39").Value = "  +9.65%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.78"
$ws.Range("E40").Value = "  +1.25%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "421.26"
$ws.Range("E41").Value = "  +4.47%  "

$ws.Range("E42").Value = "  +5.48%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.911.77"
$ws.Range("E43").Value = "  +5.04%  "

$ws.Range("E44").Value = "  +8.90%  "

$ws.Range("E45").Value = "  +0.68%  "

$ws.Range("E46").Value = "  +6.48%  "

$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "35.16"
$ws.Range("E47").Value = "  +2.17%  "

$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.999"
$ws.Range("E48").Value = "  +0.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.26"
$ws.Range("E49").Value = "  +2.92%  "

$ws.Range("E50").Value = "  +0.58%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.44"
$ws.Range("E51").Value = "  +3.85%  "
